# Edit script for "DECEMBER 21" sheet of DANIEL KATEI TIPANKO.xlsx
# Implements the rent-statement updates described in the commit diff:
#  - fills in WATER (F) and PAID (H) amounts for a number of tenants
#  - clears a stray B/F value in E25
#  - fixes the F19 number format to match its neighbours (general, style 16)
#  - reworks the expense table at rows 61-64 (new "PAID ON 11/12" /
#    "DEPOSIT REFUND" lines, adjusted amounts) and adds a helper total in L63
#  - restores the selection to H39

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DECEMBER 21")
$ws.Activate()

# ---- WATER (F) / PAID (H) updates for individual tenants ----------------

$ws.Range("H5").Value = 20000

$ws.Range("F7").Value = 300
$ws.Range("H7").Value = 20300

$ws.Range("F10").Value = 300

$ws.Range("F11").Value = 300
$ws.Range("H11").Value = 5850

$ws.Range("F12").Value = 300

$ws.Range("F13").Value = 300

$ws.Range("F14").Value = 300

$ws.Range("F15").Value = 300

$ws.Range("F16").Value = 300
$ws.Range("H16").Value = 5500

$ws.Range("F18").Value = 300

# F19 keeps its value (300) but its format should match F5:F18 (style 16,
# General number format) instead of the accounting style used by F20:F42.
$ws.Range("F18").Copy() | Out-Null
$ws.Range("F19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("F21").Value = 600
$ws.Range("H21").Value = 10000

$ws.Range("F22").Value = 750
$ws.Range("H22").Value = 11450

$ws.Range("F23").Value = 450
$ws.Range("H23").Value = 12000

# E25 previously held a stray 12000 B/F figure - clear it out.
$ws.Range("E25").ClearContents()

$ws.Range("F27").Value = 750
$ws.Range("H27").Value = 5450

$ws.Range("F28").Value = 750
$ws.Range("H28").Value = 6000

$ws.Range("F29").Value = 750

$ws.Range("F30").Value = 750

$ws.Range("F33").Value = 2550

$ws.Range("F35").Value = 1050

$ws.Range("F36").Value = 900

$ws.Range("F37").Value = 600
$ws.Range("H37").Value = 9000

$ws.Range("F38").Value = 450
$ws.Range("H38").Value = 5450

$ws.Range("F39").Value = 450
$ws.Range("H39").Formula = "=2500+2000"

$ws.Range("F40").Value = 450
$ws.Range("H40").Value = 5450

$ws.Range("F41").Value = 150
$ws.Range("H41").Value = 5000

# ---- Expense summary block (rows 61-64) ----------------------------------

# Row 61's stray helper total in J61 is no longer needed.
$ws.Range("J61").ClearContents()

# Row 62 becomes "ACUMEN SACCO" 5000 (was "DAVIS HARDWARE" 15000).
$ws.Range("A62").Value = "ACUMEN SACCO"
$ws.Range("C62").Value = 5000
$ws.Range("E62").Value = "ACUMEN SACCO"
$ws.Range("G62").Value = 5000

# Row 63 becomes "PAID ON 11/12" 148820 (was "ACUMEN SACCO" 5000), plus a
# helper figure in L63.
$ws.Range("A63").Value = "PAID ON 11/12"
$ws.Range("C63").Value = 148820
$ws.Range("E63").Value = "PAID ON 11/12"
$ws.Range("G63").Value = 148820
$ws.Range("L63").Formula = "=10000-750"

# Row 64 becomes "DEPOSIT REFUND" 9550 (was blank).
$ws.Range("A64").Value = "DEPOSIT REFUND"
$ws.Range("C64").Value = 9550
$ws.Range("E64").Value = "DEPOSIT REFUND"
$ws.Range("G64").Value = 9550

# ---- Column width / view state --------------------------------------------

# Widen column H slightly (matches the new "width=13.4..." column entry).
$ws.Columns.Item(8).ColumnWidth = 12.71

# Restore selection to H39 (matches the saved view state in the diff).
$ws.Range("H39").Select()
